$wb = $excel.ActiveWorkbook

# 1) Status text change: "Ready for handoff" -> "In Translation"
#    This shared string is referenced by the "Status"-style columns on all
#    three sheets (Overview!E2/F2, zh-cn!C2, de-de!C2), so replace it across
#    every worksheet in one pass (whole-cell match) so every reference moves
#    together.
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Cells.Replace("Ready for handoff", "In Translation", 1)
}

# 2) Narrow the report's "Status" columns (used for the archived report).
#    Overview columns E (zh-cn) and F (de-de), plus the "Status" column (C)
#    on both the zh-cn and de-de detail sheets.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
